$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (row 10), since the updated dataset has only 8 data rows (rows 2-9)
$ws.Rows(10).Delete()

# Rewrite rows 2-9 with the updated TPM-derived values
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ngfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.798983
$ws.Range("H2").Value = 11.396949
$ws.Range("I2").Value = 0.9354058228349695
$ws.Range("J2").Value = 0.9354058228349694
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3922183333333333
$ws.Range("N2").Value = 1.176655
$ws.Range("O2").Value = 0.06257714082953221
$ws.Range("P2").Value = 0.06257714082953222
$ws.Range("Q2").Value = 1.490030780621667
$ws.Range("R2").Value = 13.410277025595
$ws.Range("S2").Value = 0.05853502190830834
$ws.Range("T2").Value = 0.05853502190830835

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ngfr"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.798983
$ws.Range("H3").Value = 11.396949
$ws.Range("I3").Value = 0.9354058228349695
$ws.Range("J3").Value = 0.9354058228349694
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.787580333333334
$ws.Range("N3").Value = 17.362741
$ws.Range("O3").Value = 0.9233893441524432
$ws.Range("P3").Value = 0.9233893441524432
$ws.Range("Q3").Value = 21.98691929746767
$ws.Range("R3").Value = 197.882273677209
$ws.Range("S3").Value = 0.863743769263959
$ws.Range("T3").Value = 0.8637437692639589

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ngfr"
$ws.Range("D4").Value = "Neutrophils"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.798983
$ws.Range("H4").Value = 11.396949
$ws.Range("I4").Value = 0.9354058228349695
$ws.Range("J4").Value = 0.9354058228349694
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.034619
$ws.Range("N4").Value = 0.103857
$ws.Range("O4").Value = 0.005523347213187152
$ws.Range("P4").Value = 0.005523347213187152
$ws.Range("Q4").Value = 0.131516992477
$ws.Range("R4").Value = 1.183652932293
$ws.Range("S4").Value = 0.005166571144754563
$ws.Range("T4").Value = 0.005166571144754563

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Ntf3"
$ws.Range("C5").Value = "Ngfr"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.798983
$ws.Range("H5").Value = 11.396949
$ws.Range("I5").Value = 0.9354058228349695
$ws.Range("J5").Value = 0.9354058228349694
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05333966666666667
$ws.Range("N5").Value = 0.160019
$ws.Range("O5").Value = 0.008510167804837369
$ws.Range("P5").Value = 0.008510167804837371
$ws.Range("Q5").Value = 0.2026364868923333
$ws.Range("R5").Value = 1.823728382031
$ws.Range("S5").Value = 0.007960460517947566
$ws.Range("T5").Value = 0.007960460517947566

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntf3"
$ws.Range("C6").Value = "Ngfr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2623376666666666
$ws.Range("H6").Value = 0.787013
$ws.Range("I6").Value = 0.06459417716503056
$ws.Range("J6").Value = 0.06459417716503056
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3922183333333333
$ws.Range("N6").Value = 1.176655
$ws.Range("O6").Value = 0.06257714082953221
$ws.Range("P6").Value = 0.06257714082953222
$ws.Range("Q6").Value = 0.1028936423905555
$ws.Range("R6").Value = 0.926042781515
$ws.Range("S6").Value = 0.004042118921223871
$ws.Range("T6").Value = 0.004042118921223872

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntf3"
$ws.Range("C7").Value = "Ngfr"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2623376666666666
$ws.Range("H7").Value = 0.787013
$ws.Range("I7").Value = 0.06459417716503056
$ws.Range("J7").Value = 0.06459417716503056
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.787580333333334
$ws.Range("N7").Value = 17.362741
$ws.Range("O7").Value = 0.9233893441524432
$ws.Range("P7").Value = 0.9233893441524432
$ws.Range("Q7").Value = 1.518300320292555
$ws.Range("R7").Value = 13.664702882633
$ws.Range("S7").Value = 0.05964557488848429
$ws.Range("T7").Value = 0.05964557488848429

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Ntf3"
$ws.Range("C8").Value = "Ngfr"
$ws.Range("D8").Value = "Neutrophils"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.2623376666666666
$ws.Range("H8").Value = 0.787013
$ws.Range("I8").Value = 0.06459417716503056
$ws.Range("J8").Value = 0.06459417716503056
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.034619
$ws.Range("N8").Value = 0.103857
$ws.Range("O8").Value = 0.005523347213187152
$ws.Range("P8").Value = 0.005523347213187152
$ws.Range("Q8").Value = 0.009081867682333333
$ws.Range("R8").Value = 0.081736809141
$ws.Range("S8").Value = 0.0003567760684325887
$ws.Range("T8").Value = 0.0003567760684325887

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Ntf3"
$ws.Range("C9").Value = "Ngfr"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.2623376666666666
$ws.Range("H9").Value = 0.787013
$ws.Range("I9").Value = 0.06459417716503056
$ws.Range("J9").Value = 0.06459417716503056
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05333966666666667
$ws.Range("N9").Value = 0.160019
$ws.Range("O9").Value = 0.008510167804837369
$ws.Range("P9").Value = 0.008510167804837371
$ws.Range("Q9").Value = 0.01399300369411111
$ws.Range("R9").Value = 0.125937033247
$ws.Range("S9").Value = 0.0005497072868898043
$ws.Range("T9").Value = 0.0005497072868898044
